$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated in place
$ws.Range("A2").Value = "MSAT36"
$ws.Range("B2").Value = "Av. Rivadavia 7589"
$ws.Range("C2").Value = -34.6314059
$ws.Range("D2").Value = -58.4720002
$ws.Range("E2").Value = "https://i.ibb.co/fr7R0FD/Av-Rivadavia-7589-cara-hacia-Gral-Paz.jpg"
$ws.Range("F2").Value = "Flores"
$ws.Range("H2").Value = "CABA"
$ws.Range("I2").Value = "CABA"
$ws.Range("K2").Value = "Brasil"
$ws.Range("L2").Value = "Medianera"
$ws.Range("S2").Value = "Atacama"

# Row 3 - updated in place
$ws.Range("A3").Value = "MSAT42"
$ws.Range("B3").Value = "Av.Gaona 1770"
$ws.Range("C3").Value = -34.6110181
$ws.Range("D3").Value = -58.4530194
$ws.Range("E3").Value = "https://i.ibb.co/xFfKQp6/Av-Gaona-1770.jpg"
$ws.Range("F3").Value = "Caballito"
$ws.Range("H3").Value = "CABA"
$ws.Range("I3").Value = "CABA"
$ws.Range("K3").Value = "Brasil"
$ws.Range("L3").Value = "Medianera"
$ws.Range("S3").Value = "Atacama"

# Row 4 - new
$ws.Range("A4").Value = "MSAT43"
$ws.Range("B4").Value = "Av.Juan B.Justo 8324"
$ws.Range("C4").Value = -34.6341896
$ws.Range("D4").Value = -58.5057639
$ws.Range("E4").Value = "https://i.ibb.co/cJTHyc7/Av-Juan-B-Justo-8324.jpg"
$ws.Range("F4").Value = "Villa Luro"
$ws.Range("H4").Value = "CABA"
$ws.Range("I4").Value = "CABA"
$ws.Range("K4").Value = "Brasil"
$ws.Range("L4").Value = "Medianera"
$ws.Range("S4").Value = "Atacama"

# Row 5 - new
$ws.Range("A5").Value = "MSAT46"
$ws.Range("B5").Value = "Av. San Martín 7035 (tránsito a Provincia)"
$ws.Range("C5").Value = -34.5911444
$ws.Range("D5").Value = -58.5125877
$ws.Range("E5").Value = "https://i.ibb.co/dMzv67d/Av-San-Mart-n-7035-hacia-Provincia.jpg"
$ws.Range("F5").Value = "Devoto"
$ws.Range("H5").Value = "CABA"
$ws.Range("I5").Value = "CABA"
$ws.Range("K5").Value = "Brasil"
$ws.Range("L5").Value = "Medianera"
$ws.Range("S5").Value = "Atacama"

# Row 6 - new
$ws.Range("A6").Value = "MA63"
$ws.Range("B6").Value = "AV.TRIUNVIRATO 3700"
$ws.Range("C6").Value = -34.6117381
$ws.Range("D6").Value = -58.4233365
$ws.Range("E6").Value = "https://i.ibb.co/K6kDR7S/Rua-Jo-o-C-mara-pr-ximo-a-Av-das-Flores-Trajeto-para-Zona-Norte.jpg"
$ws.Range("H6").Value = "CABA"
$ws.Range("K6").Value = "Brasil"
$ws.Range("L6").Value = "Pantalla "
